# Fruta / hortaliza, semanal
# Insert a new weekly record for Feria Lagunitas de Puerto Montt - Pera
# at row 201, pushing the previous rows 201-225 down to 202-226.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 201 (shifts 201..225 -> 202..226)
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new weekly record
$ws.Range("A201").Value = 4
$ws.Range("B201").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C201").Value = "Los Lagos"
$ws.Range("D201").Value = 44637
$ws.Range("E201").Value = 10
$ws.Range("F201").Value = "Fruta"
$ws.Range("G201").Value = 100104
$ws.Range("H201").Value = "Frutos de pepita"
$ws.Range("I201").Value = 100104005
$ws.Range("J201").Value = "Pera"
$ws.Range("K201").Value = "Packham's Triumph"
$ws.Range("L201").Value = "Primera"
$ws.Range("M201").Value = 400
$ws.Range("N201").Value = 13000
$ws.Range("O201").Value = 14000
$ws.Range("P201").Value = 13500
$ws.Range("Q201").Value = "$/caja 15 kilos empedrada"
$ws.Range("R201").Value = "Región de O'Higgins"
$ws.Range("S201").Value = 900
$ws.Range("T201").Value = 15
